$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the current row 18 (IC1) for parts D1-D4
$ws.Rows.Item(18).Resize(4).Insert()

# Insert 1 new row before the (now shifted) row 34 (U1) for part RN1
$ws.Rows.Item(34).Insert()

# Match the row height used by the rest of the data rows (13.5pt)
$ws.Rows.Item(18).Resize(4).RowHeight = 13.5
$ws.Rows.Item(34).RowHeight = 13.5

# Fill in the new D1-D4 rows (18-21)
$ws.Range("A18").Value = "D1"
$ws.Range("B18").Value = 59.840000000000003
$ws.Range("C18").Value = -50.609999999999999
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = "top"

$ws.Range("A19").Value = "D2"
$ws.Range("B19").Value = 59.840000000000003
$ws.Range("C19").Value = -48.433332
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = "top"

$ws.Range("A20").Value = "D3"
$ws.Range("B20").Value = 59.840000000000003
$ws.Range("C20").Value = -46.256666000000003
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = "top"

$ws.Range("A21").Value = "D4"
$ws.Range("B21").Value = 59.840000000000003
$ws.Range("C21").Value = -44.079999999999998
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = "top"

# Fill in the new RN1 row (34)
$ws.Range("A34").Value = "RN1"
$ws.Range("B34").Value = 63.575800000000001
$ws.Range("C34").Value = -49.779998999999997
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = "top"
